$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "node_diffusion" worksheet between "node_history"
#    and "reserve_type".
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("reserve_type")
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "node_diffusion"
$newSheet.Move($refSheet)

$newSheet.Range("A1").Value = "node1"
$newSheet.Range("B1").Value = "node2"
$newSheet.Range("C1").Value = "diff_coeff"
$newSheet.Range("A1:C1").Style = "Normal"
$newSheet.Columns.Item(1).ColumnWidth = 9.85546875
$newSheet.Range("A9").Select()

# ---------------------------------------------------------------------
# 2. Add three new columns to the "nodes" sheet:
#      state_min           -> inserted right after in_max (before out_max)
#      is_temp              -> inserted right before residual_value
#      T_E_conversion        -> inserted right before residual_value
# ---------------------------------------------------------------------
$nodes = $wb.Worksheets.Item("nodes")

$nodes.Range("H1").EntireColumn.Insert()
$nodes.Range("H1").Value = "state_min"
$nodes.Range("H2:H4").Value = 0

$nodes.Range("M1:N1").EntireColumn.Insert()
$nodes.Range("M1").Value = "is_temp"
$nodes.Range("N1").Value = "T_E_conversion"
$nodes.Range("M2:N4").Value = 0

$nodes.Range("E17").Select()

# ---------------------------------------------------------------------
# 3. Update the active sheet / window view: "nodes" becomes the
#    selected/active sheet (previously "markets" was selected).
# ---------------------------------------------------------------------
$nodes.Activate()

$window = $excel.ActiveWindow
$window.WindowState = -4143
$window.Left = -120
$window.Top = -120
$window.Width = 18240
$window.Height = 28440
